$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in cell values for rows 2-29 (Date, Start Time, End Time) ---
$ws.Cells.Item(2,1).Value = 45901
$ws.Cells.Item(2,2).Value = 0.375
$ws.Cells.Item(2,3).Value = 0.39583333333333331
$ws.Cells.Item(3,1).Value = 45901
$ws.Cells.Item(3,2).Value = 0.39583333333333331
$ws.Cells.Item(3,3).Value = 0.41666666666666669
$ws.Cells.Item(4,1).Value = 45901
$ws.Cells.Item(4,2).Value = 0.41666666666666669
$ws.Cells.Item(4,3).Value = 0.4375
$ws.Cells.Item(5,1).Value = 45901
$ws.Cells.Item(5,2).Value = 0.4375
$ws.Cells.Item(5,3).Value = 0.45833333333333298
$ws.Cells.Item(6,1).Value = 45901
$ws.Cells.Item(6,2).Value = 0.45833333333333298
$ws.Cells.Item(6,3).Value = 0.47916666666666702
$ws.Cells.Item(7,1).Value = 45901
$ws.Cells.Item(7,2).Value = 0.5
$ws.Cells.Item(7,3).Value = 0.52083333333333337
$ws.Cells.Item(8,1).Value = 45901
$ws.Cells.Item(8,2).Value = 0.52083333333333337
$ws.Cells.Item(8,3).Value = 0.54166666666666663
$ws.Cells.Item(9,1).Value = 45901
$ws.Cells.Item(9,2).Value = 0.54166666666666663
$ws.Cells.Item(9,3).Value = 0.5625
$ws.Cells.Item(10,1).Value = 45901
$ws.Cells.Item(10,2).Value = 0.5625
$ws.Cells.Item(10,3).Value = 0.58333333333333304
$ws.Cells.Item(11,1).Value = 45901
$ws.Cells.Item(11,2).Value = 0.60416666666666663
$ws.Cells.Item(11,3).Value = 0.60416666666666663
$ws.Cells.Item(12,1).Value = 45901
$ws.Cells.Item(12,2).Value = 0.625
$ws.Cells.Item(12,3).Value = 0.625
$ws.Cells.Item(13,1).Value = 45901
$ws.Cells.Item(13,2).Value = 0.64583333333333304
$ws.Cells.Item(13,3).Value = 0.64583333333333404
$ws.Cells.Item(14,1).Value = 45901
$ws.Cells.Item(14,2).Value = 0.66666666666666696
$ws.Cells.Item(14,3).Value = 0.66666666666666696
$ws.Cells.Item(15,1).Value = 45901
$ws.Cells.Item(15,2).Value = 0.6875
$ws.Cells.Item(15,3).Value = 0.687500000000001
$ws.Cells.Item(16,1).Value = 45902
$ws.Cells.Item(16,2).Value = 0.375
$ws.Cells.Item(16,3).Value = 0.39583333333333331
$ws.Cells.Item(17,1).Value = 45902
$ws.Cells.Item(17,2).Value = 0.39583333333333331
$ws.Cells.Item(17,3).Value = 0.41666666666666669
$ws.Cells.Item(18,1).Value = 45902
$ws.Cells.Item(18,2).Value = 0.41666666666666669
$ws.Cells.Item(18,3).Value = 0.4375
$ws.Cells.Item(19,1).Value = 45902
$ws.Cells.Item(19,2).Value = 0.4375
$ws.Cells.Item(19,3).Value = 0.45833333333333298
$ws.Cells.Item(20,1).Value = 45902
$ws.Cells.Item(20,2).Value = 0.45833333333333298
$ws.Cells.Item(20,3).Value = 0.47916666666666702
$ws.Cells.Item(21,1).Value = 45902
$ws.Cells.Item(21,2).Value = 0.5
$ws.Cells.Item(21,3).Value = 0.52083333333333337
$ws.Cells.Item(22,1).Value = 45902
$ws.Cells.Item(22,2).Value = 0.52083333333333337
$ws.Cells.Item(22,3).Value = 0.54166666666666663
$ws.Cells.Item(23,1).Value = 45902
$ws.Cells.Item(23,2).Value = 0.54166666666666663
$ws.Cells.Item(23,3).Value = 0.5625
$ws.Cells.Item(24,1).Value = 45902
$ws.Cells.Item(24,2).Value = 0.5625
$ws.Cells.Item(24,3).Value = 0.58333333333333304
$ws.Cells.Item(25,1).Value = 45902
$ws.Cells.Item(25,2).Value = 0.60416666666666663
$ws.Cells.Item(25,3).Value = 0.625
$ws.Cells.Item(26,1).Value = 45902
$ws.Cells.Item(26,2).Value = 0.625
$ws.Cells.Item(26,3).Value = 0.64583333333333337
$ws.Cells.Item(27,1).Value = 45902
$ws.Cells.Item(27,2).Value = 0.64583333333333304
$ws.Cells.Item(27,3).Value = 0.66666666666666696
$ws.Cells.Item(28,1).Value = 45902
$ws.Cells.Item(28,2).Value = 0.66666666666666696
$ws.Cells.Item(28,3).Value = 0.6875
$ws.Cells.Item(29,1).Value = 45902
$ws.Cells.Item(29,2).Value = 0.6875
$ws.Cells.Item(29,3).Value = 0.70833333333333304

# --- Apply the date format (matches existing A2 style) to the newly added date cells ---
$ws.Range("A2").Copy()
$ws.Range("A6:A29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Give brand-new "wrap" time cells the same (pre-change) time format as B2:C5 ---
# so that the subsequent format change (below) merges them into one consistent style.
# (NOTE: the COM host only applies Range ops to the FIRST area of a multi-area
#  (comma) Range, so every contiguous block gets its own statement.)
$ws.Range("B2").Copy()
$ws.Range("B6:B7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B2").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B2").Copy()
$ws.Range("B16:B21").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B2").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C2").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C2").Copy()
$ws.Range("C16:C20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Switch every "wrap" time cell (old + new) from h:mm to h:mm AM/PM ---
$ws.Range("B2:B7").NumberFormat = "h:mm AM/PM"
$ws.Range("B9").NumberFormat = "h:mm AM/PM"
$ws.Range("B16:B21").NumberFormat = "h:mm AM/PM"
$ws.Range("B23").NumberFormat = "h:mm AM/PM"
$ws.Range("C2:C6").NumberFormat = "h:mm AM/PM"
$ws.Range("C16:C20").NumberFormat = "h:mm AM/PM"

# --- Brand-new "no-wrap" time cells: apply h:mm AM/PM directly (no prior style) ---
$ws.Range("B8").NumberFormat = "h:mm AM/PM"
$ws.Range("B10:B15").NumberFormat = "h:mm AM/PM"
$ws.Range("B22").NumberFormat = "h:mm AM/PM"
$ws.Range("B24:B29").NumberFormat = "h:mm AM/PM"
$ws.Range("C7:C15").NumberFormat = "h:mm AM/PM"
$ws.Range("C21:C29").NumberFormat = "h:mm AM/PM"

# --- Selection / cursor ---
[void]$ws.Range("C2").Select()
